$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2026-02-24 Tuesday", $false, $false, $false, $false,
                         $false, $true, 1, $false, "2026-02-25 Wednesday", 2)

# Update the division problems in the table. Addressed by (row, column) to
# avoid ambiguity since some cell texts repeat ("400÷2=" appears twice with
# two different replacements).
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "543÷4="
$t.Cell(1,2).Range.Text = "475÷4="
$t.Cell(1,3).Range.Text = "765÷2="
$t.Cell(1,4).Range.Text = "558÷3="
$t.Cell(1,5).Range.Text = "567÷7="

$t.Cell(5,1).Range.Text = "854÷5="
$t.Cell(5,2).Range.Text = "816÷2="
$t.Cell(5,3).Range.Text = "834÷3="
$t.Cell(5,4).Range.Text = "262÷2="
$t.Cell(5,5).Range.Text = "753÷6="

$t.Cell(9,1).Range.Text = "635÷7="
$t.Cell(9,2).Range.Text = "465÷7="
$t.Cell(9,3).Range.Text = "882÷5="
$t.Cell(9,4).Range.Text = "457÷3="
$t.Cell(9,5).Range.Text = "108÷3="

$t.Cell(13,1).Range.Text = "506÷3="
$t.Cell(13,2).Range.Text = "925÷3="
$t.Cell(13,3).Range.Text = "273÷7="
$t.Cell(13,4).Range.Text = "647÷5="
$t.Cell(13,5).Range.Text = "158÷3="

$t.Cell(17,1).Range.Text = "295÷4="
$t.Cell(17,2).Range.Text = "347÷3="
$t.Cell(17,3).Range.Text = "147÷9="
$t.Cell(17,4).Range.Text = "240÷4="
$t.Cell(17,5).Range.Text = "276÷2="
